$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 11 and 12 (data reduced from 11 rows to 9 rows of data)
$ws.Rows("11:12").Delete()

# Row 2
$ws.Range("D2").Value2 = 232.5
$ws.Range("F2").Value2 = 116.1822033898305
$ws.Range("G2").Value2 = 8
$ws.Range("H2").Value2 = "Indiana"
$ws.Range("I2").Value2 = "Boston"
$ws.Range("J2").Value2 = 0.4916741005055011
$ws.Range("K2").Value2 = 99.05225988700562
$ws.Range("L2").Value2 = 116.6006497175141
$ws.Range("M2").Value2 = 114.6872598870057
$ws.Range("N2").Value2 = 75.9145338983051
$ws.Range("O2").Value2 = 0.4540964689265535
$ws.Range("P2").Value2 = 0.5899411016949152
$ws.Range("Q2").Value2 = 0.2589987288135593
$ws.Range("R2").Value2 = 12.10983050847458
$ws.Range("S2").Value2 = 11.83005649717514
$ws.Range("T2").Value2 = 0.2088283898305084
$ws.Range("U2").Value2 = 1.01646722125836
$ws.Range("V2").Value2 = 1.029708745220725
$ws.Range("W2").Value2 = 11.07405510240506
$ws.Range("X2").Value2 = 0.5725988700564972
$ws.Range("Y2").Value2 = 39.5
$ws.Range("Z2").Value2 = 77.25
$ws.Range("AA2").Value2 = 0.5147894100714879

# Row 3
$ws.Range("D3").Value2 = 222.5
$ws.Range("F3").Value2 = 114.4194220616838
$ws.Range("G3").Value2 = 2.5
$ws.Range("H3").Value2 = "Cleveland"
$ws.Range("I3").Value2 = "Denver"
$ws.Range("J3").Value2 = 0.4752709085857183
$ws.Range("K3").Value2 = 96.3146846346207
$ws.Range("L3").Value2 = 117.9316615726591
$ws.Range("M3").Value2 = 112.8200889135871
$ws.Range("N3").Value2 = 77.62970269519312
$ws.Range("O3").Value2 = 0.3656968602389553
$ws.Range("P3").Value2 = 0.6024366490691858
$ws.Range("Q3").Value2 = 0.2685429285912754
$ws.Range("R3").Value2 = 12.47414559599889
$ws.Range("S3").Value2 = 12.62697971658794
$ws.Range("T3").Value2 = 0.2068174492914699
$ws.Range("U3").Value2 = 1.001044812438179
$ws.Range("V3").Value2 = 0.9829776545556644
$ws.Range("W3").Value2 = 10.90610648038377
$ws.Range("X3").Value2 = 0.6589330369547097
$ws.Range("Y3").Value2 = 49
$ws.Range("Z3").Value2 = 76.15
$ws.Range("AA3").Value2 = 0.4672367535391297

# Row 4
$ws.Range("D4").Value2 = 227
$ws.Range("F4").Value2 = 111.5847457627119
$ws.Range("G4").Value2 = 7
$ws.Range("H4").Value2 = "Orlando"
$ws.Range("I4").Value2 = "Detroit"
$ws.Range("J4").Value2 = 0.5214761040532365
$ws.Range("K4").Value2 = 98.66525423728817
$ws.Range("L4").Value2 = 112.3262711864406
$ws.Range("M4").Value2 = 117.5194915254237
$ws.Range("N4").Value2 = 75.50169491525423
$ws.Range("O4").Value2 = 0.3709406779661015
$ws.Range("P4").Value2 = 0.571364406779661
$ws.Range("Q4").Value2 = 0.3061949152542373
$ws.Range("R4").Value2 = 12.87033898305085
$ws.Range("S4").Value2 = 12.18135593220339
$ws.Range("T4").Value2 = 0.2298008474576271
$ws.Range("U4").Value2 = 0.976244494861871
$ws.Range("V4").Value2 = 1.016752366452033
$ws.Range("W4").Value2 = 11.19434640197002
$ws.Range("X4").Value2 = 0.3305084745762712
$ws.Range("Y4").Value2 = 28
$ws.Range("Z4").Value2 = 74.4
$ws.Range("AA4").Value2 = 0.4864970743784106

# Row 5
$ws.Range("D5").Value2 = 229
$ws.Range("F5").Value2 = 115.2368421052632
$ws.Range("G5").Value2 = 3.5
$ws.Range("H5").Value2 = "Philadelphia"
$ws.Range("I5").Value2 = "Memphis"
$ws.Range("J5").Value2 = 0.5219155844155844
$ws.Range("K5").Value2 = 98.72982456140352
$ws.Range("L5").Value2 = 116.1140350877193
$ws.Range("M5").Value2 = 112.2236842105264
$ws.Range("N5").Value2 = 76.44122807017541
$ws.Range("O5").Value2 = 0.3736578947368421
$ws.Range("P5").Value2 = 0.5827631578947369
$ws.Range("Q5").Value2 = 0.2855350877192982
$ws.Range("R5").Value2 = 11.80964912280702
$ws.Range("S5").Value2 = 12.81666666666666
$ws.Range("T5").Value2 = 0.2168114035087719
$ws.Range("U5").Value2 = 1.008196343878068
$ws.Range("V5").Value2 = 1.006548579013775
$ws.Range("W5").Value2 = 11.74619852777213
$ws.Range("X5").Value2 = 0.6403508771929824
$ws.Range("Y5").Value2 = 50
$ws.Range("Z5").Value2 = 75.7
$ws.Range("AA5").Value2 = 0.4980496117560247

# Row 6
$ws.Range("D6").Value2 = 226
$ws.Range("F6").Value2 = 114.1101694915254
$ws.Range("G6").Value2 = 6
$ws.Range("H6").Value2 = "Toronto"
$ws.Range("I6").Value2 = "NewOrleans"
$ws.Range("J6").Value2 = 0.5593220338983051
$ws.Range("K6").Value2 = 97.85847457627119
$ws.Range("L6").Value2 = 115.706779661017
$ws.Range("M6").Value2 = 114.5720338983051
$ws.Range("N6").Value2 = 77.03813559322035
$ws.Range("O6").Value2 = 0.3505508474576272
$ws.Range("P6").Value2 = 0.5695932203389831
$ws.Range("Q6").Value2 = 0.282364406779661
$ws.Range("R6").Value2 = 11.36186440677966
$ws.Range("S6").Value2 = 13.73389830508474
$ws.Range("T6").Value2 = 0.2210805084745762
$ws.Range("U6").Value2 = 0.998339190652016
$ws.Range("V6").Value2 = 0.9790972090364114
$ws.Range("W6").Value2 = 10.29313616501341
$ws.Range("X6").Value2 = 0.4915254237288136
$ws.Range("Y6").Value2 = 46
$ws.Range("Z6").Value2 = 76.15
$ws.Range("AA6").Value2 = 0.4948919009376189

# Row 7
$ws.Range("D7").Value2 = 237.5
$ws.Range("F7").Value2 = 112.7203389830509
$ws.Range("G7").Value2 = 14.5
$ws.Range("H7").Value2 = "Dallas"
$ws.Range("I7").Value2 = "SanAntonio"
$ws.Range("J7").Value2 = 0.5644736842105263
$ws.Range("K7").Value2 = 97.80063559322033
$ws.Range("L7").Value2 = 114.288884180791
$ws.Range("M7").Value2 = 119.0446610169492
$ws.Range("N7").Value2 = 75.76508474576272
$ws.Range("O7").Value2 = 0.4118829096045197
$ws.Range("P7").Value2 = 0.5778872881355933
$ws.Range("Q7").Value2 = 0.275907627118644
$ws.Range("R7").Value2 = 12.06532485875706
$ws.Range("S7").Value2 = 11.73139830508475
$ws.Range("T7").Value2 = 0.2119942796610169
$ws.Range("U7").Value2 = 0.98617969363999
$ws.Range("V7").Value2 = 1.009759102789509
$ws.Range("W7").Value2 = 11.14153430933071
$ws.Range("X7").Value2 = 0.3769774011299435
$ws.Range("Y7").Value2 = 35.5
$ws.Range("Z7").Value2 = 74.25
$ws.Range("AA7").Value2 = 0.516566190771442

# Row 8
$ws.Range("D8").Value2 = 240.5
$ws.Range("F8").Value2 = 117.730701754386
$ws.Range("G8").Value2 = 1
$ws.Range("H8").Value2 = "Utah"
$ws.Range("I8").Value2 = "OklahomaCity"
$ws.Range("J8").Value2 = 0.5691347011596788
$ws.Range("K8").Value2 = 99.96938596491228
$ws.Range("L8").Value2 = 116.8641666666667
$ws.Range("M8").Value2 = 115.9775438596491
$ws.Range("N8").Value2 = 73.76627192982455
$ws.Range("O8").Value2 = 0.3988868421052632
$ws.Range("P8").Value2 = 0.5831324561403506
$ws.Range("Q8").Value2 = 0.2619907894736842
$ws.Range("R8").Value2 = 12.03767543859649
$ws.Range("S8").Value2 = 12.53298245614035
$ws.Range("T8").Value2 = 0.2140616228070175
$ws.Range("U8").Value2 = 1.030014888489816
$ws.Range("V8").Value2 = 1.0263220926661
$ws.Range("W8").Value2 = 11.39161062614853
$ws.Range("X8").Value2 = 0.487280701754386
$ws.Range("Y8").Value2 = 23.5
$ws.Range("Z8").Value2 = 74.35
$ws.Range("AA8").Value2 = 0.4878405815256194

# Row 9
$ws.Range("D9").Value2 = 238.5
$ws.Range("F9").Value2 = 117.8188194038574
$ws.Range("G9").Value2 = 6
$ws.Range("H9").Value2 = "LALakers"
$ws.Range("I9").Value2 = "GoldenState"
$ws.Range("J9").Value2 = 0.5752118644067796
$ws.Range("K9").Value2 = 101.3956165984804
$ws.Range("L9").Value2 = 115.1713325540619
$ws.Range("M9").Value2 = 115.7277030976038
$ws.Range("N9").Value2 = 76.20055523085912
$ws.Range("O9").Value2 = 0.4113633839859732
$ws.Range("P9").Value2 = 0.5903826709526592
$ws.Range("Q9").Value2 = 0.2581227352425483
$ws.Range("R9").Value2 = 12.74040035067212
$ws.Range("S9").Value2 = 11.35952659263589
$ws.Range("T9").Value2 = 0.2032233343074226
$ws.Range("U9").Value2 = 1.030785821556058
$ws.Range("V9").Value2 = 0.9985815435608614
$ws.Range("W9").Value2 = 11.0448316658053
$ws.Range("X9").Value2 = 0.4788135593220339
$ws.Range("Y9").Value2 = 48.5
$ws.Range("Z9").Value2 = 75.65
$ws.Range("AA9").Value2 = 0.468069482260938

# Row 10
$ws.Range("D10").Value2 = 240.5
$ws.Range("F10").Value2 = 117.1677253478524
$ws.Range("G10").Value2 = 6
$ws.Range("H10").Value2 = "Sacramento"
$ws.Range("I10").Value2 = "Portland"
$ws.Range("J10").Value2 = 0.4821428571428572
$ws.Range("K10").Value2 = 98.56624319419234
$ws.Range("L10").Value2 = 118.3398215366001
$ws.Range("M10").Value2 = 117.4650786448881
$ws.Range("N10").Value2 = 76.65511191772535
$ws.Range("O10").Value2 = 0.4140003024803387
$ws.Range("P10").Value2 = 0.6044671808832427
$ws.Range("Q10").Value2 = 0.2941643980641259
$ws.Range("R10").Value2 = 12.39467634603751
$ws.Range("S10").Value2 = 11.89913793103448
$ws.Range("T10").Value2 = 0.2231632637628554
$ws.Range("U10").Value2 = 1.025089460611132
$ws.Range("V10").Value2 = 1.01463041418264
$ws.Range("W10").Value2 = 12.43665165612157
$ws.Range("X10").Value2 = 0.5220810647307925
$ws.Range("Y10").Value2 = 37
$ws.Range("Z10").Value2 = 75.5
$ws.Range("AA10").Value2 = 0.4740337871227711
